{"js": "// The worksheet is a table of unique two-digit multiplication prompts.\n// Each \"AxB=\" string occurs exactly once in the document, so a direct\n// search-and-replace per pair reproduces the diff unambiguously.\nconst replacements = [\n  [\"13\u00d732=\", \"12\u00d764=\"],\n  [\"17\u00d766=\", \"42\u00d797=\"],\n  [\"62\u00d787=\", \"29\u00d734=\"],\n  [\"44\u00d794=\", \"26\u00d748=\"],\n  [\"28\u00d746=\", \"71\u00d748=\"],\n  [\"19\u00d722=\", \"15\u00d742=\"],\n  [\"27\u00d734=\", \"78\u00d732=\"],\n  [\"52\u00d738=\", \"71\u00d719=\"],\n  [\"16\u00d756=\", \"12\u00d769=\"],\n  [\"94\u00d744=\", \"38\u00d742=\"],\n  [\"63\u00d723=\", \"24\u00d750=\"],\n  [\"62\u00d773=\", \"65\u00d714=\"],\n  [\"48\u00d730=\", \"12\u00d796=\"],\n  [\"92\u00d769=\", \"80\u00d794=\"],\n  [\"66\u00d755=\", \"65\u00d799=\"],\n  [\"89\u00d718=\", \"48\u00d740=\"],\n  [\"65\u00d723=\", \"16\u00d765=\"],\n  [\"86\u00d742=\", \"66\u00d792=\"],\n  [\"61\u00d759=\", \"85\u00d768=\"],\n  [\"61\u00d791=\", \"26\u00d729=\"],\n  [\"71\u00d723=\", \"45\u00d713=\"],\n  [\"37\u00d755=\", \"49\u00d757=\"],\n  [\"88\u00d756=\", \"14\u00d724=\"],\n  [\"25\u00d736=\", \"99\u00d790=\"],\n  [\"47\u00d790=\", \"69\u00d785=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each multiplication expression in the worksheet table is unique, so a\n# plain whole-text Find/Replace for each pair is unambiguous.\n$pairs = @(\n    @(\"13\u00d732=\", \"12\u00d764=\"),\n    @(\"17\u00d766=\", \"42\u00d797=\"),\n    @(\"62\u00d787=\", \"29\u00d734=\"),\n    @(\"44\u00d794=\", \"26\u00d748=\"),\n    @(\"28\u00d746=\", \"71\u00d748=\"),\n    @(\"19\u00d722=\", \"15\u00d742=\"),\n    @(\"27\u00d734=\", \"78\u00d732=\"),\n    @(\"52\u00d738=\", \"71\u00d719=\"),\n    @(\"16\u00d756=\", \"12\u00d769=\"),\n    @(\"94\u00d744=\", \"38\u00d742=\"),\n    @(\"63\u00d723=\", \"24\u00d750=\"),\n    @(\"62\u00d773=\", \"65\u00d714=\"),\n    @(\"48\u00d730=\", \"12\u00d796=\"),\n    @(\"92\u00d769=\", \"80\u00d794=\"),\n    @(\"66\u00d755=\", \"65\u00d799=\"),\n    @(\"89\u00d718=\", \"48\u00d740=\"),\n    @(\"65\u00d723=\", \"16\u00d765=\"),\n    @(\"86\u00d742=\", \"66\u00d792=\"),\n    @(\"61\u00d759=\", \"85\u00d768=\"),\n    @(\"61\u00d791=\", \"26\u00d729=\"),\n    @(\"71\u00d723=\", \"45\u00d713=\"),\n    @(\"37\u00d755=\", \"49\u00d757=\"),\n    @(\"88\u00d756=\", \"14\u00d724=\"),\n    @(\"25\u00d736=\", \"99\u00d790=\"),\n    @(\"47\u00d790=\", \"69\u00d785=\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
